$d = $word.ActiveDocument

# Locate the target paragraph (the one that currently holds the old
# "SmartRewards are a price stabilization..." text, identified by a
# unique substring so the lookup is resilient to index drift).
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*will get paid each 30 days on average*") {
        $targetIndex = $i
        break
    }
}
if (-not $targetIndex) { throw "Could not locate the SmartRewards paragraph" }

$target = $d.Paragraphs.Item($targetIndex)
$insertionRange = $target.Range
$insertionRange.Collapse(1)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t>SmartRewards</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> are a price stabilization mechanism and a way to encourage long term holding. Long term holders are key to the project’s success since the SmartHive treasury needs SmartCash to appreciate in value in order to fund meaningful 3rd party proposals and help grow SmartCash into a successful global crypto-currency. Beginning at block 574,100, the current SmartRewards model will change to a Decentralized Distribution. SmartRewards will then be distributed after every </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t xml:space="preserve">47,500 Blocks</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> to all eligible addresses. The payouts will begin 200 Blocks after the cycle ends and 1000 addresses will be paid every other block. SmartNodes will keep track on which addresses are eligible and which are not. You will earn SmartRewards on </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t>any</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> address for which you hold the keys (web or desktop, SmartNode included!) which holds &gt;=1000 SMART for one month and does not make </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t>any</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve">outgoing transactions during that time. Please note, most exchanges do </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:i/>
      <w:iCs/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t>not</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> pay SmartRewards to their users, holding &gt;= 1000 SMART on an exchange does not guarantee a reward</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t xml:space="preserve">. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:hAnsi="Helvetica Neue" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve">The SmartRewards will come out of the 15% block reward allocation.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t xml:space="preserve">Payment after every 47500 blocks starting at 574100. Typically, around the 25th of each month.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t xml:space="preserve">All users need to move funds into addresses holding at least 1000 SMART before the snapshot to be counted.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:textAlignment w:val="baseline"/>
    <w:rPr>
      <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
      <w:color w:val="252525"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="inherit" w:hAnsi="inherit" w:cs="Times New Roman" w:eastAsia="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="2B2B2B"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
      <w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/>
    </w:rPr>
    <w:t xml:space="preserve">If you spend ANY amount from an address, it will be ineligible for SmartRewards until the next round.</w:t>
  </w:r>
</w:p>

'@

$insertionRange.InsertXML($newXml)

# Remove the three old trailing paragraphs that described the previous
# (now superseded) payout schedule: the "Here is a small example..."
# paragraph through the "All users need to move funds into addresses
# larger than 1000 SMART..." paragraph.
$startIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Here is a small example*") {
        $startIndex = $i
        break
    }
}
if (-not $startIndex) { throw "Could not locate the start of the old payout paragraphs" }

$endIndex = $null
for ($i = $startIndex; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*addresses larger than 1000 SMART*") {
        $endIndex = $i
        break
    }
}
if (-not $endIndex) { throw "Could not locate the end of the old payout paragraphs" }

$startPara = $d.Paragraphs.Item($startIndex)
$endPara = $d.Paragraphs.Item($endIndex)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

Write-Output "Done. Paragraphs.Count=$($d.Paragraphs.Count)"
